$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet view is explicitly left-to-right (matches the saved view state).
$ws.DisplayRightToLeft = $false

$newRows = @(
    @{ r = 30; a = 'https://www.abbviecare.fr/'; b = 200 },
    @{ r = 31; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-header-publish.min.css'; b = 200 },
    @{ r = 32; a = 'https://www.abbviecare.fr/content/dam/abbvie-care-ous/fr/design-assets/banner_logo.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'; b = 200 },
    @{ r = 33; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-components.min.css'; b = 200 },
    @{ r = 34; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-footer-publish.min.js'; b = 200 },
    @{ r = 35; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-main.min.js'; b = 200 },
    @{ r = 36; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-main.min.css'; b = 200 },
    @{ r = 37; a = 'https://www.abbviecare.fr/content/dam/abbvie-care-ous/fr/design-assets/banner-desktop.jpg/_jcr_content/renditions/cq5dam.web.1280.1280.jpeg'; b = 200 },
    @{ r = 38; a = 'https://www.abbviecare.fr/etc.clientlibs/clientlibs/granite/jquery/granite/csrf.min.js'; b = 200 },
    @{ r = 39; a = 'https://www.abbviecare.fr/libs/granite/csrf/token.json'; b = 200 },
    @{ r = 40; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/assets/resources/fonts/NeueHaasUnicaPro-Regular.woff2'; b = 200 },
    @{ r = 41; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/assets/resources/fonts/NeueHaasUnicaPro-Light.woff2'; b = 200 },
    @{ r = 42; a = 'https://consent.trustarc.com/v2/notice/0ivu3f'; b = 200 },
    @{ r = 43; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/components/content/external-link-popup/v2/external-link-popup/clientlibs.min.js'; b = 200 },
    @{ r = 44; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/components/content/button-link/v2/button-link/clientlibs.min.js'; b = 200 },
    @{ r = 45; a = 'https://www.abbviecare.fr/etc.clientlibs/awcm-projects-ous/clientlibs/abbvie-care-ous/clientlibs-components.min.js'; b = 200 },
    @{ r = 46; a = 'https://consent.trustarc.com/v2/asset/ic-close.svg'; b = 200 },
    @{ r = 47; a = 'https://consent.trustarc.com/v2/asset/trustarc-logo-xs.svg'; b = 200 },
    @{ r = 48; a = 'https://consent.trustarc.com/v2/asset/latin.woff2'; b = 200 },
    @{ r = 49; a = 'https://www.abbviecare.fr/bin/public/abbvie-commons/hreflangs?resourcePath=/content/abbvie-care-ous/fr/fr/jcr:content'; b = 200 },
    @{ r = 50; a = 'https://consent.trustarc.com/v2/asset/16:19:48.8270ivu3f_AbbVieID-logo.png'; b = 200 },
    @{ r = 51; a = 'https://consent-reporting.trustarc.com/api/user-action/log?action=impression&domain=0ivu3f&behavior=implied&country=bd&language=en&rand=0.29232804318817207&session=490ba2c3-7ae8-45bb-9aae-881f8918ea78&userType=NEW'; b = 202 },
    @{ r = 52; a = 'https://consent-reporting.trustarc.com/api/user-action/bannermsg?action=views&domain=0ivu3f&behavior=implied&country=bd&language=en&rand=0.8796923997695323&session=490ba2c3-7ae8-45bb-9aae-881f8918ea78&userType=NEW'; b = 202 },
    @{ r = 53; a = 'https://consent.trustarc.com/v2/consentcategories/getnonemptyindexes?cmId=0ivu3f&referer=&fullURL=https%3A%2F%2Fwww.abbviecare.fr%2F&category='; b = 200 },
    @{ r = 54; a = 'https://consent.trustarc.com/v2/asset/ic-error.svg'; b = 200 },
    @{ r = 55; a = 'https://consent.trustarc.com/v2/asset/ic-close-white.svg'; b = 200 },
    @{ r = 56; a = 'https://www.abbviecare.fr/content/abbvie-care-ous/fr/fr/jcr:content/body/banner/banner-content/column_control/par1-100col/image_extension/item_1.coreimg.png/1663712714254-banner_logo.png'; b = 200 },
    @{ r = 57; a = 'https://www.abbviecare.fr/content/abbvie-care-ous/fr/fr/jcr:content/body/banner/bg-image/item_1.coreimg.jpg/1663712702019-banner-desktop.jpg'; b = 200 }
)

foreach ($row in $newRows) {
    $ws.Cells.Item($row.r, 1).Value = $row.a
    $ws.Cells.Item($row.r, 2).Value = $row.b
}

# Re-apply the "number stored as text" ignored-error flag over the full,
# now-larger data range (A1:B57) so the URL/status columns keep being
# treated as intentional text, same as the original A1:B29 range.
$lastRow = 1 + $newRows.Count + 28
$fullRange = $ws.Range("A1:B$lastRow")
$fullRange.Errors.Item(9).Ignore = $true
